# Regenerate save_data to use K instead of Strike# (column G "K")
# Writes the recalculated K values for rows 2-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 5
    5  = 6
    6  = 7
    7  = 5
    8  = 7
    9  = 4
    10 = 6
    11 = 9
    12 = 4
    13 = 3
    14 = 2
    15 = 8
    16 = 11
    17 = 4
    18 = 6
    19 = 8
    20 = 7
    21 = 4
    22 = 6
    23 = 10
    24 = 10
    25 = 0
    26 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
